$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.186289859955414
$ws.Range("C2").Value = 5.415197239585693
$ws.Range("D2").Value = 4.914983740478952
$ws.Range("F2").Value = 24.22138206326738
$ws.Range("G2").Value = 3.630661518984748
$ws.Range("K2").Value = 7.477755084766953
$ws.Range("M2").Value = 19.97348976505278
$ws.Range("N2").Value = 18.62662823188503
$ws.Range("O2").Value = 21.68669177339028

$ws.Range("B3").Value = 7.948335465945146
$ws.Range("C3").Value = 5.355731818822692
$ws.Range("D3").Value = 4.857600004527685
$ws.Range("F3").Value = 24.22020529099287
$ws.Range("G3").Value = 3.632356283506724
$ws.Range("K3").Value = 7.328179464958929
$ws.Range("M3").Value = 19.37736071226869
$ws.Range("N3").Value = 18.69110690134988
$ws.Range("O3").Value = 21.7339261144806

$ws.Range("B4").Value = 7.800079601986794
$ws.Range("C4").Value = 5.318415566143483
$ws.Range("D4").Value = 4.821399114178567
$ws.Range("F4").Value = 24.22559247331566
$ws.Range("G4").Value = 3.633451605803436
$ws.Range("K4").Value = 7.236007530234015
$ws.Range("M4").Value = 19.00999860558249
$ws.Range("N4").Value = 18.73244971263029
$ws.Range("O4").Value = 21.76737625797156

$ws.Range("B5").Value = 7.739221831710232
$ws.Range("C5").Value = 5.30301528217147
$ws.Range("D5").Value = 4.806410842031521
$ws.Range("F5").Value = 24.22932399207686
$ws.Range("G5").Value = 3.633911767050044
$ws.Range("K5").Value = 7.198416027914124
$ws.Range("M5").Value = 18.86021150111938
$ws.Range("N5").Value = 18.74973955706197
$ws.Range("O5").Value = 21.78212368678196

$ws.Range("B6").Value = 7.729092732754816
$ws.Range("C6").Value = 5.300446638590683
$ws.Range("D6").Value = 4.803907999747329
$ws.Range("F6").Value = 24.23003636743856
$ws.Range("G6").Value = 3.633989011832137
$ws.Range("K6").Value = 7.192173701552432
$ws.Range("M6").Value = 18.83534198650656
$ws.Range("N6").Value = 18.7526372839199
$ws.Range("O6").Value = 21.7846398329225

$ws.Range("B7").Value = 7.79926050887237
$ws.Range("C7").Value = 5.318208644951688
$ws.Range("D7").Value = 4.821197923690416
$ws.Range("F7").Value = 24.22563657883718
$ws.Range("G7").Value = 3.633457755727659
$ws.Range("K7").Value = 7.235500611648926
$ws.Range("M7").Value = 19.0079785230017
$ws.Range("N7").Value = 18.7326810965568
$ws.Range("O7").Value = 21.76757063061402

$ws.Range("B8").Value = 8.104746785986856
$ws.Range("C8").Value = 5.394865601552786
$ws.Range("D8").Value = 4.895402710687301
$ws.Range("F8").Value = 24.21970843444839
$ws.Range("G8").Value = 3.631234542449729
$ws.Range("K8").Value = 7.42628152183606
$ws.Range("M8").Value = 19.76836914013797
$ws.Range("N8").Value = 18.64849768201488
$ws.Range("O8").Value = 21.7020536615066

$ws.Range("B9").Value = 8.682734883585177
$ws.Range("C9").Value = 5.538439951289286
$ws.Range("D9").Value = 5.032946229206011
$ws.Range("F9").Value = 24.2565281185382
$ws.Range("G9").Value = 3.627306973695901
$ws.Range("K9").Value = 7.795527938816053
$ws.Range("M9").Value = 21.23850005275067
$ws.Range("N9").Value = 18.49724642951889
$ws.Range("O9").Value = 21.60896390096564

$ws.Range("B10").Value = 9.089656053759027
$ws.Range("C10").Value = 5.639331855469414
$ws.Range("D10").Value = 5.128756026609939
$ws.Range("F10").Value = 24.31300987879525
$ws.Range("G10").Value = 3.624681873631457
$ws.Range("K10").Value = 8.146764383642797
$ws.Range("M10").Value = 22.29298686047559
$ws.Range("N10").Value = 18.39444841893168
$ws.Range("O10").Value = 21.56226546308686

$ws.Range("B11").Value = 9.270032184027848
$ws.Range("C11").Value = 5.68413633824544
$ws.Range("D11").Value = 5.171126955626749
$ws.Range("F11").Value = 24.34505235033046
$ws.Range("G11").Value = 3.623543576704915
$ws.Range("K11").Value = 8.310557634028001
$ws.Range("M11").Value = 22.76476077687822
$ws.Range("N11").Value = 18.3494680174411
$ws.Range("O11").Value = 21.54575343545397

$ws.Range("B12").Value = 9.337589904351548
$ws.Range("C12").Value = 5.700937814931001
$ws.Range("D12").Value = 5.186991045335949
$ws.Range("F12").Value = 24.35809350393713
$ws.Range("G12").Value = 3.623120519650477
$ws.Range("K12").Value = 8.371522478310469
$ws.Range("M12").Value = 22.94209613387365
$ws.Range("N12").Value = 18.33268976546888
$ws.Range("O12").Value = 21.54018246552749

$ws.Range("B13").Value = 9.323074336355628
$ws.Range("C13").Value = 5.697326779805032
$ws.Range("D13").Value = 5.183582572403904
$ws.Range("F13").Value = 24.35524460254241
$ws.Range("G13").Value = 3.623211277830451
$ws.Range("K13").Value = 8.358440182984332
$ws.Range("M13").Value = 22.90396491498109
$ws.Range("N13").Value = 18.33629195210839
$ws.Range("O13").Value = 21.54135193146943

$ws.Range("B14").Value = 9.275605520230656
$ws.Range("C14").Value = 5.685521957800244
$ws.Range("D14").Value = 5.172435764171915
$ws.Range("F14").Value = 24.3461071045198
$ws.Range("G14").Value = 3.623508611636715
$ws.Range("K14").Value = 8.315594648830483
$ws.Range("M14").Value = 22.77937757238132
$ws.Range("N14").Value = 18.34808256176733
$ws.Range("O14").Value = 21.54528143736311

$ws.Range("B15").Value = 9.246430396593638
$ws.Range("C15").Value = 5.678269441257909
$ws.Range("D15").Value = 5.165584283948001
$ws.Range("F15").Value = 24.34062810590943
$ws.Range("G15").Value = 3.623691776596731
$ws.Range("K15").Value = 8.289211660901032
$ws.Range("M15").Value = 22.70288800277453
$ws.Range("N15").Value = 18.35533779409985
$ws.Range("O15").Value = 21.54777719840524

$ws.Range("B16").Value = 9.07776698852286
$ws.Range("C16").Value = 5.63638111685739
$ws.Range("D16").Value = 5.125962019040167
$ws.Range("F16").Value = 24.31104312731082
$ws.Range("G16").Value = 3.624757384633821
$ws.Range("K16").Value = 8.135913268968368
$ws.Range("M16").Value = 22.26198050951897
$ws.Range("N16").Value = 18.39742373996987
$ws.Range("O16").Value = 21.56343989618771

$ws.Range("B17").Value = 8.973035796069558
$ws.Range("C17").Value = 5.610398573375189
$ws.Range("D17").Value = 5.101339641681751
$ws.Range("F17").Value = 24.29451611614278
$ws.Range("G17").Value = 3.625425380468934
$ws.Range("K17").Value = 8.040011543830973
$ws.Range("M17").Value = 21.98933899886493
$ws.Range("N17").Value = 18.42369766437583
$ws.Range("O17").Value = 21.57426134766133

$ws.Range("B18").Value = 8.912355081687174
$ws.Range("C18").Value = 5.59535180334673
$ws.Range("D18").Value = 5.087063653202842
$ws.Range("F18").Value = 24.2856081800446
$ws.Range("G18").Value = 3.625814855736072
$ws.Range("K18").Value = 7.984178127498335
$ws.Range("M18").Value = 21.83178704199043
$ws.Range("N18").Value = 18.43897764017215
$ws.Range("O18").Value = 21.58093082609048

$ws.Range("B19").Value = 8.891735822763525
$ws.Range("C19").Value = 5.590239904441338
$ws.Range("D19").Value = 5.082210693509473
$ws.Range("F19").Value = 24.28269495892017
$ws.Range("G19").Value = 3.625947630400585
$ws.Range("K19").Value = 7.965159206997042
$ws.Range("M19").Value = 21.77832195403587
$ws.Range("N19").Value = 18.44418006271452
$ws.Range("O19").Value = 21.58326542437293

$ws.Range("B20").Value = 8.984230865338054
$ws.Range("C20").Value = 5.613175105896877
$ws.Range("D20").Value = 5.103972565600596
$ws.Range("F20").Value = 24.29621359469706
$ws.Range("G20").Value = 3.625353726918017
$ws.Range("K20").Value = 8.050290351554887
$ws.Range("M20").Value = 22.01843963730104
$ws.Range("N20").Value = 18.42088339339977
$ws.Range("O20").Value = 21.57306329131357

$ws.Range("B21").Value = 9.289569023843674
$ws.Range("C21").Value = 5.688993860977739
$ws.Range("D21").Value = 5.175714810011408
$ws.Range("F21").Value = 24.34876643306843
$ws.Range("G21").Value = 3.6234210609537
$ws.Range("K21").Value = 8.328208406365471
$ws.Range("M21").Value = 22.81600891123358
$ws.Range("N21").Value = 18.34461246956273
$ws.Range("O21").Value = 21.54410873295308

$ws.Range("B22").Value = 9.484746148742225
$ws.Range("C22").Value = 5.737580482227565
$ws.Range("D22").Value = 5.22154553799304
$ws.Range("F22").Value = 24.38839825232292
$ws.Range("G22").Value = 3.622204511379185
$ws.Range("K22").Value = 8.503653790325888
$ws.Range("M22").Value = 23.32953074665032
$ws.Range("N22").Value = 18.29624976045097
$ws.Range("O22").Value = 21.52915964975089

$ws.Range("B23").Value = 9.380997252085162
$ws.Range("C23").Value = 5.711739791366199
$ws.Range("D23").Value = 5.197183568946398
$ws.Range("F23").Value = 24.36676448794027
$ws.Range("G23").Value = 3.622849560872204
$ws.Range("K23").Value = 8.410590402018105
$ws.Range("M23").Value = 23.05621623699374
$ws.Range("N23").Value = 18.3219265003244
$ws.Range("O23").Value = 21.53677419621351

$ws.Range("B24").Value = 8.979171031384281
$ws.Range("C24").Value = 5.611920174629923
$ws.Range("D24").Value = 5.102782594780404
$ws.Range("F24").Value = 24.29544431510642
$ws.Range("G24").Value = 3.625386104554571
$ws.Range("K24").Value = 8.045645475180185
$ws.Range("M24").Value = 22.00528573782506
$ws.Range("N24").Value = 18.42215518089517
$ws.Range("O24").Value = 21.57360353702355

$ws.Range("B25").Value = 8.529160470351526
$ws.Range("C25").Value = 5.500372024385902
$ws.Range("D25").Value = 4.996630943768795
$ws.Range("F25").Value = 24.24138761580253
$ws.Range("G25").Value = 3.628323528019554
$ws.Range("K25").Value = 7.696463315071539
$ws.Range("M25").Value = 20.84442531596239
$ws.Range("N25").Value = 18.53669405370834
$ws.Range("O25").Value = 21.63034521078781
